# Recomputed NATMI ligand-receptor edge statistics with updated TPM values
# (Icam5-Itgb2 sheet) for all sending/target cluster combinations (rows 2-21).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.186041333333334
$ws.Range("H2").Value = 3.558124
$ws.Range("I2").Value = 0.3004487152423012
$ws.Range("J2").Value = 0.3004487152423012
$ws.Range("M2").Value = 0.5804443333333333
$ws.Range("N2").Value = 1.741333
$ws.Range("O2").Value = 0.002431273010151717
$ws.Range("P2").Value = 0.002435427107574628
$ws.Range("Q2").Value = 0.6884309710324445
$ws.Range("R2").Value = 6.195878739292001
$ws.Range("S2").Value = 0.0007304728523033658
$ws.Range("T2").Value = 0.0007317209455370708

# Row 3
$ws.Range("G3").Value = 1.186041333333334
$ws.Range("H3").Value = 3.558124
$ws.Range("I3").Value = 0.3004487152423012
$ws.Range("J3").Value = 0.3004487152423012
$ws.Range("O3").Value = 0.0004752041289926495
$ws.Range("P3").Value = 0.00047601606752829
$ws.Range("Q3").Value = 0.1345571799608889
$ws.Range("R3").Value = 1.211014619648
$ws.Range("S3").Value = 0.0001427744700336783
$ws.Range("T3").Value = 0.0001430184159235672

# Row 4
$ws.Range("G4").Value = 1.186041333333334
$ws.Range("H4").Value = 3.558124
$ws.Range("I4").Value = 0.3004487152423012
$ws.Range("J4").Value = 0.3004487152423012
$ws.Range("M4").Value = 136.1000366666667
$ws.Range("N4").Value = 408.30011
$ws.Range("O4").Value = 0.5700742118164518
$ws.Range("P4").Value = 0.5710482463260632
$ws.Range("Q4").Value = 161.4202689548489
$ws.Range("R4").Value = 1452.78242059364
$ws.Range("S4").Value = 0.1712780645330204
$ws.Range("T4").Value = 0.1715707119500348

# Row 5
$ws.Range("G5").Value = 1.186041333333334
$ws.Range("H5").Value = 3.558124
$ws.Range("I5").Value = 0.3004487152423012
$ws.Range("J5").Value = 0.3004487152423012
$ws.Range("M5").Value = 1.221658
$ws.Range("N5").Value = 2.443316
$ws.Range("O5").Value = 0.005117086949542552
$ws.Range("P5").Value = 0.003417220037046797
$ws.Range("Q5").Value = 1.448936883197333
$ws.Range("R5").Value = 8.693621299184001
$ws.Range("S5").Value = 0.001537422199773206
$ws.Range("T5").Value = 0.001026699369830959

# Row 6
$ws.Range("G6").Value = 1.186041333333334
$ws.Range("H6").Value = 3.558124
$ws.Range("I6").Value = 0.3004487152423012
$ws.Range("J6").Value = 0.3004487152423012
$ws.Range("M6").Value = 100.7253213333333
$ws.Range("N6").Value = 302.175964
$ws.Range("O6").Value = 0.4219022240948613
$ws.Range("P6").Value = 0.4226230904617871
$ws.Range("Q6").Value = 119.4643944146151
$ws.Range("R6").Value = 1075.179549731536
$ws.Range("S6").Value = 0.1267599811871705
$ws.Range("T6").Value = 0.1269765645609748

# Row 7
$ws.Range("I7").Value = 0.331325035675986
$ws.Range("J7").Value = 0.3313250356759861
$ws.Range("M7").Value = 0.5804443333333333
$ws.Range("N7").Value = 1.741333
$ws.Range("O7").Value = 0.002431273010151717
$ws.Range("P7").Value = 0.002435427107574628
$ws.Range("Q7").Value = 0.7591792025265556
$ws.Range("R7").Value = 6.832612822739
$ws.Range("S7").Value = 0.0008055416168265797
$ws.Range("T7").Value = 0.0008069179733034274

# Row 8
$ws.Range("I8").Value = 0.331325035675986
$ws.Range("J8").Value = 0.3313250356759861
$ws.Range("O8").Value = 0.0004752041289926495
$ws.Range("P8").Value = 0.00047601606752829
$ws.Range("S8").Value = 0.0001574470249918655
$ws.Range("T8").Value = 0.0001577160405561533

# Row 9
$ws.Range("I9").Value = 0.331325035675986
$ws.Range("J9").Value = 0.3313250356759861
$ws.Range("M9").Value = 136.1000366666667
$ws.Range("N9").Value = 408.30011
$ws.Range("O9").Value = 0.5700742118164518
$ws.Range("P9").Value = 0.5710482463260632
$ws.Range("Q9").Value = 178.0090033906811
$ws.Range("R9").Value = 1602.08103051613
$ws.Range("S9").Value = 0.1888798585680455
$ws.Range("T9").Value = 0.1892025805866922

# Row 10
$ws.Range("I10").Value = 0.331325035675986
$ws.Range("J10").Value = 0.3313250356759861
$ws.Range("M10").Value = 1.221658
$ws.Range("N10").Value = 2.443316
$ws.Range("O10").Value = 0.005117086949542552
$ws.Range("P10").Value = 0.003417220037046797
$ws.Range("Q10").Value = 1.597840297404666
$ws.Range("R10").Value = 9.587041784427999
$ws.Range("S10").Value = 0.001695419016114309
$ws.Range("T10").Value = 0.001132210550687225

# Row 11
$ws.Range("I11").Value = 0.331325035675986
$ws.Range("J11").Value = 0.3313250356759861
$ws.Range("M11").Value = 100.7253213333333
$ws.Range("N11").Value = 302.175964
$ws.Range("O11").Value = 0.4219022240948613
$ws.Range("P11").Value = 0.4226230904617871
$ws.Range("Q11").Value = 131.7414345057569
$ws.Range("R11").Value = 1185.672910551812
$ws.Range("S11").Value = 0.1397867694500078
$ws.Range("T11").Value = 0.1400256105247471

# Row 12
$ws.Range("G12").Value = 0.6500023333333333
$ws.Range("H12").Value = 1.950007
$ws.Range("I12").Value = 0.1646589882374797
$ws.Range("J12").Value = 0.1646589882374797
$ws.Range("M12").Value = 0.5804443333333333
$ws.Range("N12").Value = 1.741333
$ws.Range("O12").Value = 0.002431273010151717
$ws.Range("P12").Value = 0.002435427107574628
$ws.Range("Q12").Value = 0.3772901710367778
$ws.Range("R12").Value = 3.395611539331
$ws.Range("S12").Value = 0.0004003309539806734
$ws.Range("T12").Value = 0.0004010149634593698

# Row 13
$ws.Range("G13").Value = 0.6500023333333333
$ws.Range("H13").Value = 1.950007
$ws.Range("I13").Value = 0.1646589882374797
$ws.Range("J13").Value = 0.1646589882374797
$ws.Range("O13").Value = 0.0004752041289926495
$ws.Range("P13").Value = 0.00047601606752829
$ws.Range("Q13").Value = 0.07374319805155555
$ws.Range("R13").Value = 0.663688782464
$ws.Range("S13").Value = 0.00007824663108620243
$ws.Range("T13").Value = 0.00007838032406399204

# Row 14
$ws.Range("G14").Value = 0.6500023333333333
$ws.Range("H14").Value = 1.950007
$ws.Range("I14").Value = 0.1646589882374797
$ws.Range("J14").Value = 0.1646589882374797
$ws.Range("M14").Value = 136.1000366666667
$ws.Range("N14").Value = 408.30011
$ws.Range("O14").Value = 0.5700742118164518
$ws.Range("P14").Value = 0.5710482463260632
$ws.Range("Q14").Value = 88.46534140008556
$ws.Range("R14").Value = 796.1880726007701
$ws.Range("S14").Value = 0.09386784293797562
$ws.Range("T14").Value = 0.09402822647483662

# Row 15
$ws.Range("G15").Value = 0.6500023333333333
$ws.Range("H15").Value = 1.950007
$ws.Range("I15").Value = 0.1646589882374797
$ws.Range("J15").Value = 0.1646589882374797
$ws.Range("M15").Value = 1.221658
$ws.Range("N15").Value = 2.443316
$ws.Range("O15").Value = 0.005117086949542552
$ws.Range("P15").Value = 0.003417220037046797
$ws.Range("Q15").Value = 0.7940805505353333
$ws.Range("R15").Value = 4.764483303212
$ws.Range("S15").Value = 0.0008425743598348877
$ws.Range("T15").Value = 0.0005626759938849684

# Row 16
$ws.Range("G16").Value = 0.6500023333333333
$ws.Range("H16").Value = 1.950007
$ws.Range("I16").Value = 0.1646589882374797
$ws.Range("J16").Value = 0.1646589882374797
$ws.Range("M16").Value = 100.7253213333333
$ws.Range("N16").Value = 302.175964
$ws.Range("O16").Value = 0.4219022240948613
$ws.Range("P16").Value = 0.4226230904617871
$ws.Range("Q16").Value = 65.47169389241645
$ws.Range("R16").Value = 589.245245031748
$ws.Range("S16").Value = 0.06946999335460227
$ws.Range("T16").Value = 0.0695886904812347

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.8035953333333333
$ws.Range("H17").Value = 2.410786
$ws.Range("I17").Value = 0.2035672608442332
$ws.Range("J17").Value = 0.2035672608442332
$ws.Range("M17").Value = 0.5804443333333333
$ws.Range("N17").Value = 1.741333
$ws.Range("O17").Value = 0.002431273010151717
$ws.Range("P17").Value = 0.002435427107574628
$ws.Range("Q17").Value = 0.4664423575264445
$ws.Range("R17").Value = 4.197981217738
$ws.Range("S17").Value = 0.0004949275870410986
$ws.Range("T17").Value = 0.0004957732252747607

# Row 18
$ws.Range("E18").Value = 3
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 0.8035953333333333
$ws.Range("H18").Value = 2.410786
$ws.Range("I18").Value = 0.2035672608442332
$ws.Range("J18").Value = 0.2035672608442332
$ws.Range("O18").Value = 0.0004752041289926495
$ws.Range("P18").Value = 0.00047601606752829
$ws.Range("Q18").Value = 0.09116842629688889
$ws.Range("R18").Value = 0.8205158366719999
$ws.Range("S18").Value = 0.0000967360028809033
$ws.Range("T18").Value = 0.00009690128698457753

# Row 19
$ws.Range("E19").Value = 3
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 0.8035953333333333
$ws.Range("H19").Value = 2.410786
$ws.Range("I19").Value = 0.2035672608442332
$ws.Range("J19").Value = 0.2035672608442332
$ws.Range("M19").Value = 136.1000366666667
$ws.Range("N19").Value = 408.30011
$ws.Range("O19").Value = 0.5700742118164518
$ws.Range("P19").Value = 0.5710482463260632
$ws.Range("Q19").Value = 109.3693543318289
$ws.Range("R19").Value = 984.3241889864599
$ws.Range("S19").Value = 0.1160484457774103
$ws.Range("T19").Value = 0.1162467273144996

# Row 20
$ws.Range("E20").Value = 3
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 0.8035953333333333
$ws.Range("H20").Value = 2.410786
$ws.Range("I20").Value = 0.2035672608442332
$ws.Range("J20").Value = 0.2035672608442332
$ws.Range("M20").Value = 1.221658
$ws.Range("N20").Value = 2.443316
$ws.Range("O20").Value = 0.005117086949542552
$ws.Range("P20").Value = 0.003417220037046797
$ws.Range("Q20").Value = 0.9817186677293332
$ws.Range("R20").Value = 5.890312006375999
$ws.Range("S20").Value = 0.00104167137382015
$ws.Range("T20").Value = 0.0006956341226436454

# Row 21
$ws.Range("E21").Value = 3
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 0.8035953333333333
$ws.Range("H21").Value = 2.410786
$ws.Range("I21").Value = 0.2035672608442332
$ws.Range("J21").Value = 0.2035672608442332
$ws.Range("M21").Value = 100.7253213333333
$ws.Range("N21").Value = 302.175964
$ws.Range("O21").Value = 0.4219022240948613
$ws.Range("P21").Value = 0.4226230904617871
$ws.Range("Q21").Value = 80.94239817196711
$ws.Range("R21").Value = 728.481583547704
$ws.Range("S21").Value = 0.08588548010308074
$ws.Range("T21").Value = 0.08603222489483056
